$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the two new columns (I, J), copying the header
# style from H1 so the new header cells match the existing bold/border/
# centered style used by the other header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for the new I0 / IF columns, one triple per data row:
# (row, I value, J value)
$data = @(
    @(2, 8, 9),
    @(3, 7, 8),
    @(4, 7, 8),
    @(5, 5, 5),
    @(6, 8, 8),
    @(7, 7, 7),
    @(8, 7, 8),
    @(9, 8, 8),
    @(10, 7, 7),
    @(11, 8, 8),
    @(12, 8, 9),
    @(13, 8, 8),
    @(14, 8, 9),
    @(15, 7, 7),
    @(16, 8, 9),
    @(17, 6, 7),
    @(18, 7, 8),
    @(19, 8, 8),
    @(20, 6, 7),
    @(21, 8, 11),
    @(22, 7, 8),
    @(23, 6, 7),
    @(24, 10, 10),
    @(25, 8, 8),
    @(26, 8, 8),
    @(27, 6, 7),
    @(28, 6, 6),
    @(29, 7, 7),
    @(30, 6, 6),
    @(31, 5, 7),
    @(32, 7, 7),
    @(33, 8, 8),
    @(34, 7, 8),
    @(35, 8, 8),
    @(36, 8, 8),
    @(37, 8, 8),
    @(38, 7, 8),
    @(39, 8, 8),
    @(40, 8, 8),
    @(41, 8, 9),
    @(42, 4, 4),
    @(43, 3, 3),
    @(44, 3, 3)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
